$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:J1").EntireRow.Insert()

$ws.Range("A1").Value = "Date: "
$ws.Range("B1").Value = (Get-Date -Year 2022 -Month 9 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B1").NumberFormat = "m/d/yyyy"

$ws.Range("G5").Select()
